$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C4 value from 11 to 10
$ws.Range("C4").Value = 10

# Update D4, D6, D8 text from "非三角形" to "不构成三角形"
$ws.Range("D4").Value = "不构成三角形"
$ws.Range("D6").Value = "不构成三角形"
$ws.Range("D8").Value = "不构成三角形"
